# Replace the counterbalanced-character placeholders [CS1]/[CS2] with the
# actual stimulus names used in this replication (BERGMITE / PALPITOAD).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("instructions")

$cells = @("A2", "A3")
foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $text = $rng.Value()
    $text = $text -replace [regex]::Escape("[CS1]"), "BERGMITE"
    $text = $text -replace [regex]::Escape("[CS2]"), "PALPITOAD"
    $rng.Value = $text
}
